# Word COM-interop script implementing the documented diff:
#   - extend the "const keyword ... means that " paragraph with more text
#   - add several new paragraphs describing the mouse.cpp / mouse.h changes
#   - add a short "1.J " paragraph plus a couple of blank paragraphs
#   - relocate the (hidden) _GoBack bookmark to sit after the very last
#     edit, matching the fact that Word itself stamps _GoBack at the
#     location of the most recent edit.

$d = $word.ActiveDocument

# --- locate the anchor paragraph ("...Also, using const keyword means that ") ---
$find = $d.Content
$find.Find.ClearFormatting()
$null = $find.Find.Execute("const keyword")
$find.Collapse(0)  # wdCollapseEnd

# Insert the new text right after "const keyword" (i.e. ahead of the
# pre-existing " means that " run + the hidden _GoBack bookmark that used
# to sit there).
$find.InsertAfter(" means that ")
$find.Collapse(0)
$find.InsertAfter("there is confidence that the underground object in that function won")
$find.Collapse(0)
$find.InsertAfter([char]0x2019)
$find.Collapse(0)
$find.InsertAfter("t change the state of the object. For ")
$find.Collapse(0)
$find.InsertAfter("example,")
$find.Collapse(0)
$find.InsertAfter(" properties values cannot be changed.")
$find.Collapse(0)

# The paragraph that currently holds all this text.
$p1 = $find.Paragraphs(1)

# The original " means that " run (plus the hidden _GoBack bookmark that
# sat right before it) is now stranded at the end of the paragraph -
# remove it since its text has effectively been retyped above.
$tail = $p1.Range.Duplicate()
$tail.Collapse(0)              # end of paragraph (after the final period)
$null = $tail.MoveEnd(1, -1)   # step back before the paragraph mark
$null = $tail.MoveStart(1, -12) # length of " means that "
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$tail.Delete()

# --- new blank paragraph ---
$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()

# --- "In the mouse header ..." paragraph ---
$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$r3 = $p3.Range
$r3.Collapse(1) # wdCollapseStart
$r3.InsertAfter("In the mouse header, the prototype for that function needs to be changed")
$r3.Collapse(0)
$r3.InsertAfter(" to pass a constant ")
$r3.Collapse(0)
$r3.InsertAfter("underground ")
$r3.Collapse(0)
$r3.InsertAfter("reference ")
$r3.Collapse(0)
$r3.InsertAfter("object ")
$r3.Collapse(0)
$r3.InsertAfter("and ")
$r3.Collapse(0)
$r3.InsertAfter("make sure the function is const")
$r3.Collapse(0)
$r3.InsertAfter("ant.")
$r3.Collapse(0)

# --- new blank paragraph ---
$r3.InsertParagraphAfter()
$p4 = $p3.Next()

# --- "Another change would be ..." paragraph ---
$p4.Range.InsertParagraphAfter()
$p5 = $p4.Next()
$r5 = $p5.Range
$r5.Collapse(1)
$r5.InsertAfter("Another change would be that is ")
$r5.Collapse(0)
$r5.InsertAfter("that the inner function is ")
$r5.Collapse(0)
$r5.InsertAfter("at position would be ")
$r5.Collapse(0)
$r5.InsertAfter("constant.")
$r5.Collapse(0)

# The _GoBack bookmark follows the most recent edit in real Word; move it
# here, right after "constant." (end of this paragraph's text).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$gobackRange = $d.Range($r5.End, $r5.End)
$d.Bookmarks.Add("_GoBack", $gobackRange)

# --- "1.J " paragraph ---
$r5.InsertParagraphAfter()
$p6 = $p5.Next()
$r6 = $p6.Range
$r6.Collapse(1)
$r6.InsertAfter("1.J ")
$r6.Collapse(0)

# --- two trailing blank paragraphs ---
$r6.InsertParagraphAfter()
$p7 = $p6.Next()
$p7.Range.InsertParagraphAfter()
